$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: "Find The Lucky Integer in an Array"
$ws.Cells.Item(31, 1).Value = 1394
$ws.Cells.Item(31, 2).Value = "Find Lucky Integer in an Array"
$ws.Cells.Item(31, 3).Value = "Frequency Table/Dictionary"
$ws.Cells.Item(31, 4).Value = "make a freq table[arr.Length] because the lucky number must be equal to or less than the Length. Foreach to count freq, reverse loop to return the highest"
$ws.Cells.Item(31, 5).Value = "have a hashmap, count the kvp, extract the highest key from the hashmap."

# Update the view: scroll/top-left cell and selection to match the new row
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 9
$ws.Range("E31").Select()
